# Add 2022-Q4 data.
#
# The source workbook has sheets: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q2.
# This change inserts a brand-new "2022-Q4" sheet (with its own fund-holdings
# table) right after "总计" and before "2022-Q3", and adds a matching summary
# row at the top of the "总计" table. All the other quarterly sheets keep
# their existing name + data untouched; they just shift one tab to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row (same layout as the other quarterly sheets).
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B-G hold text (fund codes / names / figures rendered as strings),
# column A is the running index, column H is a numeric rank.
$q4.Range("B2:G6").NumberFormat = "@"

# Row 2 - 000593
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "000593"
$q4.Range("C2").Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$q4.Range("D2").Value = "2.30"
$q4.Range("E2").Value = "93.71"
$q4.Range("F2").Value = "9.04"
$q4.Range("G2").Value = "0.2079"
$q4.Range("H2").Value = 1

# Row 3 - 005676
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "005676"
$q4.Range("C3").Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$q4.Range("D3").Value = "2.30"
$q4.Range("E3").Value = "93.71"
$q4.Range("F3").Value = "9.04"
$q4.Range("G3").Value = "0.2079"
$q4.Range("H3").Value = 1

# Row 4 - 118002
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "118002"
$q4.Range("C4").Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$q4.Range("D4").Value = "2.30"
$q4.Range("E4").Value = "93.71"
$q4.Range("F4").Value = "9.04"
$q4.Range("G4").Value = "0.2079"
$q4.Range("H4").Value = 1

# Row 5 - 513080
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "513080"
$q4.Range("C5").Value = "华安法国CAC40ETF（QDII）"
$q4.Range("D5").Value = "0.64"
$q4.Range("E5").Value = "93.56"
$q4.Range("F5").Value = "10.63"
$q4.Range("G5").Value = "0.0680"
$q4.Range("H5").Value = 1

# Row 6 - 006282
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "006282"
$q4.Range("C6").Value = "上投摩根欧洲动力策略股票（QDII）"
$q4.Range("D6").Value = "0.43"
$q4.Range("E6").Value = "92.90"
$q4.Range("F6").Value = "3.18"
$q4.Range("G6").Value = "0.0137"
$q4.Range("H6").Value = 3

# Match the index column's (A2:A6) look to the rest of the workbook by
# copying the formatting that the other quarterly sheets use for it.
$q4.Range("A2:A6").Copy() | Out-Null
$total.Range("A2").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122) | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null

# Re-apply values clobbered by the format-only paste above (PasteSpecial of
# formats shouldn't touch values, but keep things explicit/robust).
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1
$q4.Range("A4").Value = 2
$q4.Range("A5").Value = 3
$q4.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new top row for 2022-Q4
#    and shift the previous rows down by one.
# ---------------------------------------------------------------------
$total.Range("B6").Value = "2021-Q2"
$total.Range("C6").Value = 8
$total.Range("D6").Value = 2.5

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 8
$total.Range("D5").Value = 0.75

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.58

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 5
$total.Range("D3").Value = 0.61

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.71

$total.Range("A5").Copy() | Out-Null
$total.Range("A6").PasteSpecial(-4122) | Out-Null
$total.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 3. Restore "总计" as the active sheet (it was active before the edit).
# ---------------------------------------------------------------------
$total.Activate()
